# The workbook's single sheet contains daily price records for
# "Hortaliza, Macroferia Regional de Talca - Zanahoria".
# A new daily record is inserted as row 219 (shifting every following
# row down by one), and the used range grows from A1:R293 to A1:R294.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 219; this shifts rows 219:293
# down to 220:294 and carries the existing formatting (incl. the date
# number format on column D) down with them.
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A219").Value = 5
$ws.Range("B219").Value = "Macroferia Regional de Talca"
$ws.Range("C219").Value = "Maule"
$ws.Range("D219").Value = 44627
$ws.Range("E219").Value = 7
$ws.Range("F219").Value = 100114013
$ws.Range("G219").Value = "Zanahoria"
$ws.Range("H219").Value = "Sin especificar"
$ws.Range("I219").Value = "Primera"
$ws.Range("J219").Value = 300
$ws.Range("K219").Value = 7000
$ws.Range("L219").Value = 7000
$ws.Range("M219").Value = 7000
$ws.Range("N219").Value = "$/saco 20 kilos"
$ws.Range("O219").Value = "Región de Ñuble"
$ws.Range("P219").Value = 350
$ws.Range("Q219").Value = 20
$ws.Range("R219").Value = "Hortaliza"
